$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 currently holds (via shared strings):
#   A1 = title template            (unchanged)
#   B1 = comment.is_locked_lbl...  -> drop (replaced by old D1's order_by text)
#   C1 = comment.is_enabled_lbl... -> drop (replaced by old E1's rem text)
#   D1 = comment.order_by          -> shifts left into B1
#   E1 = comment.rem               -> shifts left into C1
# A new "tenant_id" column template takes the now-freed last slot (D1),
# and the now-unused trailing cell (E1) is removed entirely.

$orderByText = $ws.Range("D1").Value()
$remText = $ws.Range("E1").Value()
$tenantIdText = '<%=comment.tenant_id_lbl%><%selectList.tenant_id = data.findAllTenant.map((item) => item.lbl)%><%_dataValidation_({ sqref: `${ _col }2:${ _col }${ _lastRow }`, formula1: `"${ selectList.tenant_id.join(",") }"` })%>'

$ws.Range("B1").Value = $orderByText
$ws.Range("C1").Value = $remText
$ws.Range("D1").Value = $tenantIdText
$ws.Range("E1").ClearContents()
